$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 45224
$ws.Range("C3").Value = 45224
$ws.Range("C4").Value = 45224
$ws.Range("C5").Value = 45224
